# Auto-generated edit script applying numeric updates described in the diff.
# Updates the "想去人数" (F) and one "最低票价" (G) column values across the
# four worksheets: 展览, 演出, 本地生活, 全部类型.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 554
$ws.Range("F4").Value = 6067
$ws.Range("F7").Value = 1586
$ws.Range("F8").Value = 20
$ws.Range("F9").Value = 34
$ws.Range("F10").Value = 703
$ws.Range("F11").Value = 1848
$ws.Range("F12").Value = 1848
$ws.Range("F13").Value = 9
$ws.Range("F14").Value = 1684
$ws.Range("F15").Value = 585
$ws.Range("F16").Value = 222
$ws.Range("F17").Value = 663
$ws.Range("F18").Value = 4792
$ws.Range("G18").Value = 59.9
$ws.Range("F19").Value = 129
$ws.Range("F21").Value = 683
$ws.Range("F22").Value = 3372
$ws.Range("F23").Value = 838
$ws.Range("F25").Value = 60
$ws.Range("F27").Value = 2376
$ws.Range("F29").Value = 354
$ws.Range("F31").Value = 14
$ws.Range("F33").Value = 1256
$ws.Range("F35").Value = 37
$ws.Range("F36").Value = 13
$ws.Range("F38").Value = 1334
$ws.Range("F39").Value = 1310
$ws.Range("F40").Value = 89

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F10").Value = 94
$ws.Range("F11").Value = 18
$ws.Range("F14").Value = 100
$ws.Range("F19").Value = 133
$ws.Range("F20").Value = 318
$ws.Range("F21").Value = 247
$ws.Range("F22").Value = 509

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F3").Value = 804
$ws.Range("F4").Value = 225
$ws.Range("F5").Value = 325

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 554
$ws.Range("F5").Value = 804
$ws.Range("F6").Value = 225
$ws.Range("F7").Value = 6067
$ws.Range("F16").Value = 94
$ws.Range("F17").Value = 1586
$ws.Range("F18").Value = 18
$ws.Range("F19").Value = 20
$ws.Range("F20").Value = 34
$ws.Range("F21").Value = 1848
$ws.Range("F22").Value = 9
$ws.Range("F23").Value = 1684
$ws.Range("F24").Value = 100
$ws.Range("F25").Value = 585
$ws.Range("F26").Value = 222
$ws.Range("F27").Value = 663
$ws.Range("F28").Value = 4792
$ws.Range("F30").Value = 683
$ws.Range("F31").Value = 3372
$ws.Range("F33").Value = 60
$ws.Range("F36").Value = 2376
$ws.Range("F37").Value = 354
$ws.Range("F40").Value = 1256
$ws.Range("F41").Value = 133
$ws.Range("F42").Value = 247
$ws.Range("F43").Value = 509
$ws.Range("F45").Value = 37
$ws.Range("F46").Value = 13
$ws.Range("F48").Value = 1334
$ws.Range("F50").Value = 89
